# nzgrapher navigation update:
#  - remove "Part 3.5: Dotplots and Informal Confidence Intervals" row entirely
#    (row 10), shifting all subsequent rows up by one
#  - swap the now-outdated NZGrapher video id H15YmE99iKI for the new one
#    mZI-W7w__r0 everywhere it is referenced (Part 3 row, and the two NCEA
#    standard cross-reference rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole "Part 3.5" row (row 10) - Excel shifts everything below
# up and re-writes the relative formulas automatically, just like a manual
# right-click > Delete on the row header.
$ws.Rows(10).Delete()

# Replace every occurrence of the old video id with the new one.
$ws.Cells.Replace("H15YmE99iKI", "mZI-W7w__r0")

# Match the author's last selection in the saved file.
$ws.Range("B14").Select()
